# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet that
#    shows it (Overview!E2:F2, zh-cn!C2, de-de!C2 - all four cells share the
#    same underlying string, so updating all four keeps the shared string
#    table collapsed onto a single entry, same as the source edit).
# 2) Narrow the "Latest Handoff Datetime" / duplicate status-width columns:
#    Overview columns E & F, and column C on the zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- status text ---------------------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# --- column widths --------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
